# Error Calculations and Plots
# The "H 72" record (row 2) is removed from the missing-data sheet; all
# rows below it shift up by one, so the used range shrinks from
# A1:F63 to A1:F62.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(2).Delete()
